$d = $word.ActiveDocument

# 1) Update the convenio date in the intro paragraph: "4 de marzo" -> "7 de marzo"
$d.Content.Find.Execute(
    "suscrito con fecha  4 de marzo de 2022",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "suscrito con fecha  7 de marzo de 2022",
    2
) | Out-Null

# 2) Update the signature-block date: "En Puertollano a  4  de marzo  2022" -> "...7..."
$d.Content.Find.Execute(
    "En Puertollano a  4  de marzo  2022",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "En Puertollano a  7  de marzo  2022",
    2
) | Out-Null

# 3) Remove the now-stale student row ("Jiménez Coello Daniel") from the
#    students table - annexes are deleted before re-associating students
#    to companies, so this leftover duplicate row goes away.
foreach ($tbl in $d.Tables) {
    for ($r = $tbl.Rows.Count; $r -ge 1; $r--) {
        $row = $tbl.Rows.Item($r)
        if ($row.Cells.Count -ge 1 -and $row.Cells.Item(1).Range.Text -like "*Jim*nez Coello Daniel*") {
            $row.Delete()
        }
    }
}
